# Sync attendance_reports: swap order of "Recorded By" names in column G
# Every cell in column G whose value is exactly "dnasr281@gmail.com, System"
# is changed to "System, dnasr281@gmail.com". All other cells (including the
# "Recorded By" header and plain "dnasr281@gmail.com" values) stay untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldValue = "dnasr281@gmail.com, System"
$newValue = "System, dnasr281@gmail.com"

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G
    if ($cell.Value2 -eq $oldValue) {
        $cell.Value = $newValue
    }
}
